$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (AD1:AF1) - match the header style used by the
# existing headers (bold font, thin border, centered/top aligned).
$ws.Range("AD1:AF1").Font.Bold = $true
$ws.Range("AD1:AF1").HorizontalAlignment = -4108
$ws.Range("AD1:AF1").VerticalAlignment = -4160
$ws.Range("AD1:AF1").Borders.LineStyle = 1

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Team record repeated for every player row (2-45)
$ws.Range("AD2:AD45").Value = 84
$ws.Range("AE2:AE45").Value = 78
$ws.Range("AF2:AF45").Value = 0
